# Auto-generated COM script applying the scheduled-runner recalculation diff
# to the per-leve profit columns (H:N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Cells.Item(2, 8).Value = 245.82353
$ws.Cells.Item(2, 9).Value = 245.82353
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 245.82353
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -132.82353
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(129, 8).Value = 952.6279
$ws.Cells.Item(129, 9).Value = 499.03705
$ws.Cells.Item(129, 10).Value = 1718.0625
$ws.Cells.Item(129, 11).Value = 1497.11115
$ws.Cells.Item(129, 12).Value = 5154.1875
$ws.Cells.Item(129, 13).Value = 3502.88885
$ws.Cells.Item(129, 14).Value = -15154.1875

$ws.Cells.Item(137, 8).Value = 2143892.2
$ws.Cells.Item(137, 9).Value = 981361.5
$ws.Cells.Item(137, 10).Value = 5264369
$ws.Cells.Item(137, 11).Value = 2944084.5
$ws.Cells.Item(137, 12).Value = 15793107
$ws.Cells.Item(137, 13).Value = -2941534.5
$ws.Cells.Item(137, 14).Value = -15798207

$ws.Cells.Item(138, 8).Value = 227349.92
$ws.Cells.Item(138, 9).Value = 944.5227
$ws.Cells.Item(138, 10).Value = 891472.4
$ws.Cells.Item(138, 11).Value = 2833.5681
$ws.Cells.Item(138, 12).Value = 2674417.2
$ws.Cells.Item(138, 13).Value = 2306.4319
$ws.Cells.Item(138, 14).Value = -2684697.2

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Cells.Item(32, 8).Value = 2828542.8
$ws.Cells.Item(32, 9).Value = 5242.104
$ws.Cells.Item(32, 10).Value = 20944722
$ws.Cells.Item(32, 11).Value = 5242.104
$ws.Cells.Item(32, 12).Value = 20944722
$ws.Cells.Item(32, 13).Value = -4955.104
$ws.Cells.Item(32, 14).Value = -20945296

$ws.Cells.Item(61, 8).Value = 1139.8269
$ws.Cells.Item(61, 9).Value = 1129.42
$ws.Cells.Item(61, 10).Value = 1400
$ws.Cells.Item(61, 11).Value = 1129.42
$ws.Cells.Item(61, 12).Value = 1400
$ws.Cells.Item(61, 13).Value = -917.4200000000001
$ws.Cells.Item(61, 14).Value = -1824

$ws.Cells.Item(74, 8).Value = 825.3182
$ws.Cells.Item(74, 9).Value = 766.2679000000001
$ws.Cells.Item(74, 10).Value = 1156
$ws.Cells.Item(74, 11).Value = 766.2679000000001
$ws.Cells.Item(74, 12).Value = 1156
$ws.Cells.Item(74, 13).Value = 107.7320999999999
$ws.Cells.Item(74, 14).Value = -2904

$ws.Cells.Item(77, 8).Value = 825.3182
$ws.Cells.Item(77, 9).Value = 766.2679000000001
$ws.Cells.Item(77, 10).Value = 1156
$ws.Cells.Item(77, 11).Value = 3831.3395
$ws.Cells.Item(77, 12).Value = 5780
$ws.Cells.Item(77, 13).Value = 536.6605
$ws.Cells.Item(77, 14).Value = -14516

$ws.Cells.Item(122, 8).Value = 2154.6155
$ws.Cells.Item(122, 9).Value = 1933.6364
$ws.Cells.Item(122, 11).Value = 5800.9092
$ws.Cells.Item(122, 13).Value = -3350.9092

$ws.Cells.Item(132, 8).Value = 89696.75999999999
$ws.Cells.Item(132, 9).Value = 129316.74
$ws.Cells.Item(132, 10).Value = 3853.4443
$ws.Cells.Item(132, 11).Value = 387950.22
$ws.Cells.Item(132, 12).Value = 11560.3329
$ws.Cells.Item(132, 13).Value = -385420.22
$ws.Cells.Item(132, 14).Value = -16620.3329

$ws.Cells.Item(136, 8).Value = 1139.8269
$ws.Cells.Item(136, 9).Value = 1129.42
$ws.Cells.Item(136, 10).Value = 1400
$ws.Cells.Item(136, 11).Value = 3388.26
$ws.Cells.Item(136, 12).Value = 4200
$ws.Cells.Item(136, 13).Value = -838.2600000000002
$ws.Cells.Item(136, 14).Value = -9300

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Cells.Item(20, 8).Value = 11908247
$ws.Cells.Item(20, 9).Value = 19612548
$ws.Cells.Item(20, 10).Value = 1598.8182
$ws.Cells.Item(20, 11).Value = 19612548
$ws.Cells.Item(20, 12).Value = 1598.8182
$ws.Cells.Item(20, 13).Value = -19612301
$ws.Cells.Item(20, 14).Value = -2092.8182

$ws.Cells.Item(134, 8).Value = 69626.17999999999
$ws.Cells.Item(134, 9).Value = 108454.43
$ws.Cells.Item(134, 10).Value = 1676.75
$ws.Cells.Item(134, 11).Value = 325363.29
$ws.Cells.Item(134, 12).Value = 5030.25
$ws.Cells.Item(134, 13).Value = -322828.29
$ws.Cells.Item(134, 14).Value = -10100.25

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Cells.Item(31, 8).Value = 2039.2941
$ws.Cells.Item(31, 9).Value = 1548.3077
$ws.Cells.Item(31, 10).Value = 3635
$ws.Cells.Item(31, 11).Value = 1548.3077
$ws.Cells.Item(31, 12).Value = 3635
$ws.Cells.Item(31, 13).Value = -1253.3077
$ws.Cells.Item(31, 14).Value = -4225

$ws.Cells.Item(34, 8).Value = 2039.2941
$ws.Cells.Item(34, 9).Value = 1548.3077
$ws.Cells.Item(34, 10).Value = 3635
$ws.Cells.Item(34, 11).Value = 1548.3077
$ws.Cells.Item(34, 12).Value = 3635
$ws.Cells.Item(34, 13).Value = -1346.3077
$ws.Cells.Item(34, 14).Value = -4039

$ws.Cells.Item(58, 8).Value = 1342.2142
$ws.Cells.Item(58, 9).Value = 1342.2142
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 1342.2142
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -1139.2142
$ws.Cells.Item(58, 14).ClearContents()

$ws.Cells.Item(99, 8).Value = 1477.7727
$ws.Cells.Item(99, 9).Value = 1334.0667
$ws.Cells.Item(99, 10).Value = 1785.7142
$ws.Cells.Item(99, 11).Value = 1334.0667
$ws.Cells.Item(99, 12).Value = 1785.7142
$ws.Cells.Item(99, 13).Value = 163.9332999999999
$ws.Cells.Item(99, 14).Value = -4781.7142

$ws.Cells.Item(107, 8).Value = 1512.7368
$ws.Cells.Item(107, 9).Value = 1571.375
$ws.Cells.Item(107, 10).Value = 1200
$ws.Cells.Item(107, 11).Value = 1571.375
$ws.Cells.Item(107, 12).Value = 1200
$ws.Cells.Item(107, 13).Value = 348.625
$ws.Cells.Item(107, 14).Value = -5040

$ws.Cells.Item(126, 8).Value = 1477.7727
$ws.Cells.Item(126, 9).Value = 1334.0667
$ws.Cells.Item(126, 10).Value = 1785.7142
$ws.Cells.Item(126, 11).Value = 4002.2001
$ws.Cells.Item(126, 12).Value = 5357.142599999999
$ws.Cells.Item(126, 13).Value = -1532.2001
$ws.Cells.Item(126, 14).Value = -10297.1426

$ws.Cells.Item(132, 8).Value = 1358.95
$ws.Cells.Item(132, 9).Value = 1198.0278
$ws.Cells.Item(132, 10).Value = 2807.25
$ws.Cells.Item(132, 11).Value = 3594.0834
$ws.Cells.Item(132, 12).Value = 8421.75
$ws.Cells.Item(132, 13).Value = -1064.0834
$ws.Cells.Item(132, 14).Value = -13481.75

$ws.Cells.Item(136, 8).Value = 1342.2142
$ws.Cells.Item(136, 9).Value = 1342.2142
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 4026.6426
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -1476.6426
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Cells.Item(5, 8).Value = 635363.7
$ws.Cells.Item(5, 9).Value = 447.6316
$ws.Cells.Item(5, 10).Value = 6667066.5
$ws.Cells.Item(5, 11).Value = 1342.8948
$ws.Cells.Item(5, 12).Value = 20001199.5
$ws.Cells.Item(5, 13).Value = -1230.8948
$ws.Cells.Item(5, 14).Value = -20001423.5

$ws.Cells.Item(107, 8).Value = 821.4583
$ws.Cells.Item(107, 9).Value = 431.625
$ws.Cells.Item(107, 10).Value = 1601.125
$ws.Cells.Item(107, 11).Value = 1294.875
$ws.Cells.Item(107, 12).Value = 4803.375
$ws.Cells.Item(107, 13).Value = 625.125
$ws.Cells.Item(107, 14).Value = -8643.375

$ws.Cells.Item(122, 8).Value = 43134.94
$ws.Cells.Item(122, 9).Value = 413.5
$ws.Cells.Item(122, 10).Value = 50089.594
$ws.Cells.Item(122, 11).Value = 3721.5
$ws.Cells.Item(122, 12).Value = 450806.346
$ws.Cells.Item(122, 13).Value = -1271.5
$ws.Cells.Item(122, 14).Value = -455706.346

$ws.Cells.Item(135, 8).Value = 635363.7
$ws.Cells.Item(135, 9).Value = 447.6316
$ws.Cells.Item(135, 10).Value = 6667066.5
$ws.Cells.Item(135, 11).Value = 4028.6844
$ws.Cells.Item(135, 12).Value = 60003598.5
$ws.Cells.Item(135, 13).Value = -1493.6844
$ws.Cells.Item(135, 14).Value = -60008668.5

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Cells.Item(102, 8).Value = 1607.814
$ws.Cells.Item(102, 9).Value = 1133.5
$ws.Cells.Item(102, 10).Value = 2333.2354
$ws.Cells.Item(102, 11).Value = 1133.5
$ws.Cells.Item(102, 12).Value = 2333.2354
$ws.Cells.Item(102, 13).Value = 488.5
$ws.Cells.Item(102, 14).Value = -5577.2354

$ws.Cells.Item(126, 8).Value = 5329.346
$ws.Cells.Item(126, 9).Value = 2236.077
$ws.Cells.Item(126, 10).Value = 8422.615
$ws.Cells.Item(126, 11).Value = 6708.231000000001
$ws.Cells.Item(126, 12).Value = 25267.845
$ws.Cells.Item(126, 13).Value = -4238.231000000001
$ws.Cells.Item(126, 14).Value = -30207.845

$ws.Cells.Item(132, 8).Value = 1534.9836
$ws.Cells.Item(132, 9).Value = 1136.0889
$ws.Cells.Item(132, 10).Value = 2656.875
$ws.Cells.Item(132, 11).Value = 3408.2667
$ws.Cells.Item(132, 12).Value = 7970.625
$ws.Cells.Item(132, 13).Value = -878.2667000000001
$ws.Cells.Item(132, 14).Value = -13030.625

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Cells.Item(40, 8).Value = 1286.4
$ws.Cells.Item(40, 9).Value = 1266.125
$ws.Cells.Item(40, 11).Value = 1266.125
$ws.Cells.Item(40, 13).Value = -1130.125

$ws.Cells.Item(95, 8).Value = 28795.143
$ws.Cells.Item(95, 10).Value = 28795.143
$ws.Cells.Item(95, 12).Value = 28795.143
$ws.Cells.Item(95, 14).Value = -34287.143

$ws.Cells.Item(118, 8).Value = 32896
$ws.Cells.Item(118, 10).Value = 32896
$ws.Cells.Item(118, 12).Value = 32896
$ws.Cells.Item(118, 14).Value = -36210

$ws.Cells.Item(122, 8).Value = 2661.2104
$ws.Cells.Item(122, 9).Value = 2088.7778
$ws.Cells.Item(122, 10).Value = 3176.4
$ws.Cells.Item(122, 11).Value = 6266.3334
$ws.Cells.Item(122, 12).Value = 9529.200000000001
$ws.Cells.Item(122, 13).Value = -3816.3334
$ws.Cells.Item(122, 14).Value = -14429.2

$ws.Cells.Item(132, 8).Value = 2754
$ws.Cells.Item(132, 9).Value = 2979.2173
$ws.Cells.Item(132, 10).Value = 2322.3333
$ws.Cells.Item(132, 11).Value = 8937.651899999999
$ws.Cells.Item(132, 12).Value = 6966.999899999999
$ws.Cells.Item(132, 13).Value = -6407.651899999999
$ws.Cells.Item(132, 14).Value = -12026.9999

$ws.Cells.Item(136, 8).Value = 1712.5366
$ws.Cells.Item(136, 9).Value = 1454.2972
$ws.Cells.Item(136, 11).Value = 4362.8916
$ws.Cells.Item(136, 13).Value = -1812.8916

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Cells.Item(126, 8).Value = 2917.6667
$ws.Cells.Item(126, 9).Value = 1835.3334
$ws.Cells.Item(126, 10).Value = 4000
$ws.Cells.Item(126, 11).Value = 5506.0002
$ws.Cells.Item(126, 12).Value = 12000
$ws.Cells.Item(126, 13).Value = -3036.0002
$ws.Cells.Item(126, 14).Value = -16940

$ws.Cells.Item(132, 8).Value = 2266.9424
$ws.Cells.Item(132, 9).Value = 2923.8064
$ws.Cells.Item(132, 10).Value = 1297.2858
$ws.Cells.Item(132, 11).Value = 8771.4192
$ws.Cells.Item(132, 12).Value = 3891.8574
$ws.Cells.Item(132, 13).Value = -6241.4192
$ws.Cells.Item(132, 14).Value = -8951.857400000001
